$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 42632.880810185183
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = "Buy"
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 10509
$ws.Range("F3").Value = 1499
$ws.Range("G3").Value = 61
$ws.Range("H3").Value = 36
$ws.Range("I3").Value = 78
$ws.Range("J3").Value = 21
$ws.Range("K3").Value = 13147
$ws.Range("L3").Value = 278
$ws.Range("M3").Value = 168
$ws.Range("N3").Value = 22
$ws.Range("O3").Value = 6
$ws.Range("P3").Value = "Noun"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.56
$ws.Range("S3").Value = 0.10150000000000001
$ws.Range("S3").NumberFormat = $ws.Range("S2").NumberFormat
$ws.Range("T3").Value = -0.93
$ws.Range("U3").Value = 2.32
$ws.Range("V3").Value = "N/A"
$ws.Range("W3").Value = 0
